# Kotte2014 workbook update — add new transport/exchange/biomass reactions
# and rework the FBP efflux row into an EC_Biomass reaction.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kotte2014")

# Cell writes are ordered to reproduce the author's original shared-string
# insertion sequence (new unique strings are appended to the shared string
# table in first-use order).

# New row 10: PEPt2r
$ws.Range("A10").Value = "PEPt2r"

# Row 4: FBP efflux reaction repurposed as EC_Biomass reaction
$ws.Range("A4").Value = "EC_Biomass"

# New row 10 (equation)
$ws.Range("C10").Value = "pep[c] <==> pep[e]"

# New row 11: PEPex
$ws.Range("A11").Value = "PEPex"
$ws.Range("C11").Value = "pep[e] <==>"

# New row 12: ENZ1ex
$ws.Range("A12").Value = "ENZ1ex"

# New row 13: ENZex
$ws.Range("A13").Value = "ENZex"

# Row 8: ENZt2r row's equation gains the enz1[e] compartment term
$ws.Range("C8").Value = "enz1[c] <==> enz1[e]"

# Row 9: ENZtr row's equation gains the enz[e] compartment term
$ws.Range("C9").Value = "enz[c] <==> enz[e]"

# Row 4 (equation)
$ws.Range("C4").Value = "fdp[c] ---> bm[c]"

# New row 14: FDex
$ws.Range("A14").Value = "FDex"
$ws.Range("C14").Value = "bm[e] <==>"

# New row 12 (equation)
$ws.Range("C12").Value = "enz1[e] <==> "

# New row 13 (equation)
$ws.Range("C13").Value = "enz[e] <==>"

# Row 7: FDPt2r row's equation changes from fdp[e] <==> to bm[c] ---> bm[e]
$ws.Range("C7").Value = "bm[c] ---> bm[e]"

# Fill in the remaining numeric columns for the new reaction rows (10 & 11;
# rows 12-14 have no E/I/J data, matching the target layout).
$ws.Range("E10").Value = 0
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1

$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1

# --- Column widths for the now-visible A & B columns ---
$ws.Columns.Item(1).ColumnWidth = 10.877604166666666
$ws.Columns.Item(2).ColumnWidth = 8.307291666666666

# --- Selection moves from C8:C9 to A2:A4 ---
$ws.Range("A2:A4").Select() | Out-Null

Write-Output "Kotte2014 sheet updated: EC_Biomass + PEP/ENZ1/ENZ/FD exchange reactions added"
